$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: Coin (B), Link (C), Price (D), Volume(1h) (E)
# The Price column value is prefixed with an apostrophe so Excel keeps it
# as text (matching the original inlineStr cells) instead of coercing it
# into a number, which would drop formatting such as trailing zeros or
# the dotted thousands separators used for some of the larger prices.
$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "'92.935.19", "  -5.61%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "'3.350.57", "  -4.77%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'1.00", "  +0.04%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'231.97", "  -8.87%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'629.23", "  -5.95%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'1.35", "  -9.47%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.386", "  -10.23%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'1.00", "  +0.11%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.931", "  -11.80%  ")
    ,@("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "'3.348.37", "  -4.93%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.194", "  -7.99%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'40.29", "  -12.59%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'5.94", "  -4.33%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "'93.054.46", "  -5.42%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "'3.969.52", "  -4.65%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.0000243", "  -6.91%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'7.94", "  -11.96%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "'3.344.44", "  -4.58%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'16.80", "  -10.78%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'10.93", "  -8.10%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'490.78", "  -6.61%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.449", "  -15.22%  ")
    ,@("SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "'3.12", "  -9.44%  ")
    ,@("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "'0.0000185", "  -10.46%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "'6.25", "  -8.23%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'89.21", "  -9.03%  ")
    ,@("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'11.42", "  -10.17%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'11.28", "  -9.68%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'1.00", "  +0.09%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "'2.64", "  -9.74%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.131", "  -10.27%  ")
    ,@("Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "'0.995", "  -0.62%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "'0.171", "  -10.86%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'28.40", "  -8.32%  ")
    ,@("PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "'0.522", "  -11.93%  ")
    ,@("RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "'7.44", "  -7.72%  ")
    ,@("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "'520.90", "  -2.27%  ")
    ,@("USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "'1.00", "  +0.03%  ")
    ,@("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "'1.39", "  -9.95%  ")
    ,@("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "'0.147", "  -5.89%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "'0.871", "  -5.27%  ")
    ,@("WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "'24.03", "  -1.74%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'1.67", "  -6.71%  ")
    ,@("MantraDAO", "https://coinranking.com/coin/cTdD8lD-6+mantradao-om", "'3.56", "  -2.91%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'5.45", "  -6.36%  ")
    ,@("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "'2.14", "  -3.73%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.0393", "  -9.33%  ")
    ,@("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "'52.07", "  -5.99%  ")
    ,@("dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "'3.10", "  -4.92%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "'7.87", "  -9.60%  ")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $r++
}

Write-Host "Updated cryptos list"